$d = $word.ActiveDocument

# 1. Fix the typo in the run that only contains "tod" (it sits between the
#    "...automatizando " run and the "o o processo..." run). Retyping just
#    that run's text as "dot" reproduces the author's "todo" -> "doto" slip
#    and lets Word's editor naturally merge the three runs into one, while
#    keeping the formatting/rsid of the surrounding (identical) runs - this
#    matches the committed change exactly, instead of a blanket
#    find-and-replace that would also hit the many other "tod..." words
#    ("todos", "todas", ...) elsewhere in the paragraph.
$rng = $d.Content.Duplicate
[void]$rng.Find.Execute("automatizando ")
$todStart = $rng.End
$todRange = $d.Range($todStart, $todStart + 3)
if ($todRange.Text -ne "tod") {
    throw "Unexpected text where the 'tod' run was expected: [" + $todRange.Text + "]"
}
$todRange.Text = "dot"

# 2. Move the "_GoBack" bookmark from right after the "O site" paragraph to
#    the empty paragraph that follows the "...compra da estufa." paragraph.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*compra da estufa.*") {
        $target = $p.Next().Range.Duplicate
    }
}
if ($target -eq $null) {
    throw "Could not locate the destination paragraph for the _GoBack bookmark"
}
$target.Collapse(1)
$d.Bookmarks.Add("_GoBack", $target)
